# attendance.xlsx: add a "sex" column between full_name and timestamp,
# and populate the sheet with the full attendance log (rows 2-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C ("sex"), pushing the existing "timestamp"
# column from C to D.
$ws.Columns.Item(3).Insert()

# Header row.
$ws.Cells.Item(1, 3).Value = "sex"

# Attendance rows: event_date, full_name, sex, timestamp.
$rows = @(
    @("2026-01-04", "Darren Nathanael Budiman", "M", "2026-01-03T11:38:59"),
    @("2026-01-04", "Samantha Gracia",           "F", "2026-01-03T11:41:02"),
    @("2026-01-03", "Shiro Budiman",              "M", "2026-01-03T11:41:35"),
    @("2026-01-04", "Ricky Adikurnia",            "M", "2026-01-03T14:15:59"),
    @("2026-01-04", "Siat Cynthia",                "F", "2026-01-03T14:15:59"),
    @("2026-01-03", "Leonardo Kwan",               "M", "2026-01-03T14:23:09"),
    @("2026-01-03", "Darren Nathanael Budiman",    "M", "2026-01-03T14:23:28"),
    @("2026-01-03", "Wilson Thiesman",             "M", "2026-01-03T19:40:07"),
    @("2026-01-03", "Haydee Judith Manuella",      "F", "2026-01-03T19:40:43"),
    @("2026-01-03", "Darren Nathanael Budiman",    "M", "2026-01-03T19:40:43"),
    @("2026-01-03", "Marcello Pardede",            "M", "2026-01-03T19:50:20"),
    @("2026-01-04", "Darren Nathanael Budiman",    "M", "2026-01-03T19:52:34"),
    @("2026-01-04", "Leonardo Kwan",               "M", "2026-01-03T19:52:34"),
    @("2026-01-03", "Richard",                     "M", "2026-01-03T19:52:50")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    # event_date: force text so Excel doesn't coerce "2026-01-0x" into a
    # date serial number - prefix with an apostrophe, then reset the
    # cell style so the quote-prefix formatting doesn't stick around.
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = "'" + $data[0]
    $cell.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
